$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "/Array_1::Vector"
$ws.Range("C1").Value = "/Array_2::Vector{Int}"

$ws.Range("C1").Select()
